$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all target cells first so values such as percentages,
# dates/timestamps, and negative decimals are stored as literal text (matching
# the original inlineStr/text cell type) rather than being auto-converted by Excel.
$targetCells = @(
  "E2",
  "O2",
  "E3",
  "O3",
  "E4",
  "H4",
  "J4",
  "E5",
  "K5",
  "N5",
  "E6",
  "O6",
  "E7",
  "J7",
  "O7",
  "E8",
  "H8",
  "N8",
  "O8",
  "E9",
  "O9",
  "E10",
  "K10",
  "N10",
  "O10",
  "E11",
  "H11",
  "N11",
  "O11",
  "E12",
  "O12",
  "E13",
  "J13",
  "K13",
  "O13",
  "E14",
  "E15",
  "H15",
  "N15",
  "O15",
  "E16",
  "H16",
  "L16",
  "O16",
  "E17",
  "N17",
  "E18",
  "N18",
  "O18",
  "E19",
  "H19",
  "K19",
  "O19",
  "E20",
  "K20",
  "N20",
  "O20",
  "E21",
  "J21",
  "K21",
  "E22",
  "O22",
  "E23",
  "H23",
  "I23",
  "K23",
  "L23",
  "E24",
  "O24",
  "E25",
  "H25",
  "K25",
  "E26",
  "J26",
  "K26",
  "E27",
  "E28",
  "J28",
  "N28",
  "O28",
  "E29",
  "N29",
  "E30",
  "J30",
  "E31",
  "J31",
  "L31",
  "E32",
  "K32",
  "O32",
  "E33",
  "J33",
  "N33",
  "O33",
  "E34",
  "H34",
  "M34",
  "O34",
  "E35",
  "J35",
  "N35",
  "O35",
  "E36",
  "E37",
  "H37",
  "J37",
  "L37",
  "O37",
  "E38",
  "H38",
  "K38",
  "N38",
  "O38",
  "E39",
  "L39",
  "E40",
  "O40",
  "E41",
  "L41",
  "O41",
  "E42",
  "E43",
  "H43",
  "N43",
  "O43",
  "E44",
  "H44",
  "L44",
  "E45",
  "J45",
  "K45",
  "N45",
  "O45",
  "E46",
  "H46",
  "J46",
  "K46",
  "N46"
)
foreach ($addr in $targetCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (meteocat automatic daily refresh).
$ws.Range('E2').Value = '2026-02-23 05:18:40'
$ws.Range('O2').Value = '3.6 °C'
$ws.Range('E3').Value = '2026-02-23 05:18:42'
$ws.Range('O3').Value = '1.5 °C'
$ws.Range('E4').Value = '2026-02-23 05:18:45'
$ws.Range('H4').Value = '91%'
$ws.Range('J4').Value = '1026.0 hPa'
$ws.Range('E5').Value = '2026-02-23 05:18:48'
$ws.Range('K5').Value = '-0.1 MJ/m2'
$ws.Range('N5').Value = '1.8 °C 4:41 TU'
$ws.Range('E6').Value = '2026-02-23 05:18:50'
$ws.Range('O6').Value = '9.4 °C'
$ws.Range('E7').Value = '2026-02-23 05:18:53'
$ws.Range('J7').Value = '1025.1 hPa'
$ws.Range('O7').Value = '11.9 °C'
$ws.Range('E8').Value = '2026-02-23 05:18:55'
$ws.Range('H8').Value = '55%'
$ws.Range('N8').Value = '11.9 °C 4:59 TU'
$ws.Range('O8').Value = '13.1 °C'
$ws.Range('E9').Value = '2026-02-23 05:18:58'
$ws.Range('O9').Value = '7.0 °C'
$ws.Range('E10').Value = '2026-02-23 05:19:01'
$ws.Range('K10').Value = '-0.1 MJ/m2'
$ws.Range('N10').Value = '3.3 °C 4:34 TU'
$ws.Range('O10').Value = '4.3 °C'
$ws.Range('E11').Value = '2026-02-23 05:19:03'
$ws.Range('H11').Value = '93%'
$ws.Range('N11').Value = '1.9 °C 4:49 TU'
$ws.Range('O11').Value = '2.9 °C'
$ws.Range('E12').Value = '2026-02-23 05:19:06'
$ws.Range('O12').Value = '5.7 °C'
$ws.Range('E13').Value = '2026-02-23 05:19:08'
$ws.Range('J13').Value = '1032.1 hPa'
$ws.Range('K13').Value = '-0.1 MJ/m2'
$ws.Range('O13').Value = '-0.9 °C'
$ws.Range('E14').Value = '2026-02-23 05:19:11'
$ws.Range('E15').Value = '2026-02-23 05:19:14'
$ws.Range('H15').Value = '88%'
$ws.Range('N15').Value = '4.7 °C 4:38 TU'
$ws.Range('O15').Value = '6.8 °C'
$ws.Range('E16').Value = '2026-02-23 05:19:16'
$ws.Range('H16').Value = '20%'
$ws.Range('L16').Value = '35.3 km/h - 211º 4:35 TU'
$ws.Range('O16').Value = '2.9 °C'
$ws.Range('E17').Value = '2026-02-23 05:19:19'
$ws.Range('N17').Value = '6.5 °C 4:59 TU'
$ws.Range('E18').Value = '2026-02-23 05:19:22'
$ws.Range('N18').Value = '1.9 °C 4:34 TU'
$ws.Range('O18').Value = '3.1 °C'
$ws.Range('E19').Value = '2026-02-23 05:19:24'
$ws.Range('H19').Value = '44%'
$ws.Range('K19').Value = '-0.1 MJ/m2'
$ws.Range('O19').Value = '9.9 °C'
$ws.Range('E20').Value = '2026-02-23 05:19:27'
$ws.Range('K20').Value = '-0.1 MJ/m2'
$ws.Range('N20').Value = '1.8 °C 4:54 TU'
$ws.Range('O20').Value = '3.2 °C'
$ws.Range('E21').Value = '2026-02-23 05:19:30'
$ws.Range('J21').Value = '1029.2 hPa'
$ws.Range('K21').Value = '-0.1 MJ/m2'
$ws.Range('E22').Value = '2026-02-23 05:19:32'
$ws.Range('O22').Value = '2.0 °C'
$ws.Range('E23').Value = '2026-02-23 05:19:34'
$ws.Range('H23').Value = '28%'
$ws.Range('I23').Value = '0.5 mm'
$ws.Range('K23').Value = '-0.1 MJ/m2'
$ws.Range('L23').Value = '34.9 km/h - 330º 4:53 TU'
$ws.Range('E24').Value = '2026-02-23 05:19:38'
$ws.Range('O24').Value = '2.6 °C'
$ws.Range('E25').Value = '2026-02-23 05:19:41'
$ws.Range('H25').Value = '31%'
$ws.Range('K25').Value = '-0.1 MJ/m2'
$ws.Range('E26').Value = '2026-02-23 05:19:43'
$ws.Range('J26').Value = '1026.7 hPa'
$ws.Range('K26').Value = '-0.1 MJ/m2'
$ws.Range('E27').Value = '2026-02-23 05:19:46'
$ws.Range('E28').Value = '2026-02-23 05:19:49'
$ws.Range('J28').Value = '1027.4 hPa'
$ws.Range('N28').Value = '2.3 °C 4:35 TU'
$ws.Range('O28').Value = '3.8 °C'
$ws.Range('E29').Value = '2026-02-23 05:19:51'
$ws.Range('N29').Value = '3.5 °C 4:59 TU'
$ws.Range('E30').Value = '2026-02-23 05:19:54'
$ws.Range('J30').Value = '1025.6 hPa'
$ws.Range('E31').Value = '2026-02-23 05:19:57'
$ws.Range('J31').Value = '1024.5 hPa'
$ws.Range('L31').Value = '53.6 km/h - 331º 4:59 TU'
$ws.Range('E32').Value = '2026-02-23 05:19:59'
$ws.Range('K32').Value = '-0.1 MJ/m2'
$ws.Range('O32').Value = '1.6 °C'
$ws.Range('E33').Value = '2026-02-23 05:20:02'
$ws.Range('J33').Value = '1029.6 hPa'
$ws.Range('N33').Value = '1.4 °C 4:40 TU'
$ws.Range('O33').Value = '2.6 °C'
$ws.Range('E34').Value = '2026-02-23 05:20:05'
$ws.Range('H34').Value = '48%'
$ws.Range('M34').Value = '4.9 °C 4:33 TU'
$ws.Range('O34').Value = '2.1 °C'
$ws.Range('E35').Value = '2026-02-23 05:20:07'
$ws.Range('J35').Value = '1026.6 hPa'
$ws.Range('N35').Value = '8.0 °C 4:44 TU'
$ws.Range('O35').Value = '10.4 °C'
$ws.Range('E36').Value = '2026-02-23 05:20:10'
$ws.Range('E37').Value = '2026-02-23 05:20:13'
$ws.Range('H37').Value = '80%'
$ws.Range('J37').Value = '1029.6 hPa'
$ws.Range('L37').Value = '12.6 km/h - 85º 4:34 TU'
$ws.Range('O37').Value = '3.8 °C'
$ws.Range('E38').Value = '2026-02-23 05:20:15'
$ws.Range('H38').Value = '71%'
$ws.Range('K38').Value = '-0.1 MJ/m2'
$ws.Range('N38').Value = '3.4 °C 4:53 TU'
$ws.Range('O38').Value = '6.7 °C'
$ws.Range('E39').Value = '2026-02-23 05:20:18'
$ws.Range('L39').Value = '24.8 km/h - 326º 4:43 TU'
$ws.Range('E40').Value = '2026-02-23 05:20:21'
$ws.Range('O40').Value = '2.1 °C'
$ws.Range('E41').Value = '2026-02-23 05:20:23'
$ws.Range('L41').Value = '9.4 km/h - 15º 4:49 TU'
$ws.Range('O41').Value = '7.3 °C'
$ws.Range('E42').Value = '2026-02-23 05:20:26'
$ws.Range('E43').Value = '2026-02-23 05:20:28'
$ws.Range('H43').Value = '94%'
$ws.Range('N43').Value = '2.4 °C 4:46 TU'
$ws.Range('O43').Value = '4.1 °C'
$ws.Range('E44').Value = '2026-02-23 05:20:31'
$ws.Range('H44').Value = '38%'
$ws.Range('L44').Value = '51.5 km/h - 46º 4:57 TU'
$ws.Range('E45').Value = '2026-02-23 05:20:33'
$ws.Range('J45').Value = '1030.6 hPa'
$ws.Range('K45').Value = '-0.1 MJ/m2'
$ws.Range('N45').Value = '2.0 °C 4:43 TU'
$ws.Range('O45').Value = '3.8 °C'
$ws.Range('E46').Value = '2026-02-23 05:20:36'
$ws.Range('H46').Value = '99%'
$ws.Range('J46').Value = '1027.4 hPa'
$ws.Range('K46').Value = '-0.1 MJ/m2'
$ws.Range('N46').Value = '0.8 °C 4:58 TU'
